$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.387.37"
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").Value = "2.355.54"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.28%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "542.90"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.95%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "135.58"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("E7").Value = "  +0.60%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.561"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +4.80%  "
$ws.Range("E9").Value = "  +0.32%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "5.66"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +6.81%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.152"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.51%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.357"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.47%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "23.92"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "2.775.18"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").Value = "58.333.85"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "2.339.72"
$ws.Range("E17").Value = "  +0.36%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "10.76"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.97%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "333.59"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.12%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.27"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.36%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.72"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.89%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "62.80"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("B24").Value = "Kaspa"
$ws.Range("C24").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.168"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.87%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "8.46"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.40"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +5.03%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.76"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "170.53"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0738"
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.13"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "18.44"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("B33").Value = "SuiNetwork"
$ws.Range("C33").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.02"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +12.52%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.25"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +5.91%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.25"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.65"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +4.08%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "39.17"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "143.30"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -3.76%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.378"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.65"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "289.19"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.87%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0942"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "19.22"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.50%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0503"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.566"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0220"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("B49").Value = "Polygon"
$ws.Range("C49").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.385"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "17.58"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "11.07"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.37%  "
